$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Biotechnology tickers: remove AEZS
$ws.Range("B7").Value = "RUBY, SLGL, OCUP"

# Internet Retail tickers: add OCG and duplicate RMBL
$ws.Range("B8").Value = "OCG, RMBL, RMBL"

# Last industry row: change from "Banks - Regional" / GNTY to "Drug Manufacturers - Specialty Generic" / ADMS
$ws.Range("A16").Value = "Drug Manufacturers - Specialty Generic "
$ws.Range("B16").Value = "ADMS"
